# Reverse the order of all worksheet tabs (2020-Q4 .. 总计  ->  总计 .. 2020-Q4)
# while keeping each sheet's own data/name together, and keep "2020-Q4"
# (the sheet that was active before the edit) the active/selected tab
# afterwards, now that it sits in the last position.

$wb = $excel.ActiveWorkbook

$count = $wb.Worksheets.Count
for ($i = 1; $i -le $count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Move($wb.Worksheets.Item(1))
}

$wb.Worksheets.Item("2020-Q4").Activate()
